$wb = $excel.ActiveWorkbook

# ---- Sheet: snapshot ----
$ws1 = $wb.Worksheets.Item("snapshot")

$ws1.Cells.Item(2,1).Value = "АВТ"
$ws1.Cells.Item(2,2).Value = "Автомобилист"
$ws1.Cells.Item(2,3).Value = "avtomobilist"
$ws1.Cells.Item(2,4).Value = "Да Коста Стефан"
$ws1.Cells.Item(2,5).Value = "'77"
$ws1.Cells.Item(2,6).Value = "нападающий"
$ws1.Cells.Item(2,7).Value = "'22668"
$ws1.Cells.Item(2,8).Value = "1369_АВТ_дакостастефан"
$ws1.Cells.Item(2,9).Value = "injured_active"
$ws1.Cells.Item(2,10).Value = "https://www.khl.ru/clubs/avtomobilist/team/"
$ws1.Cells.Item(2,11).Value = "2025-12-13T03:01:21.677037+00:00"

$ws1.Cells.Item(3,1).Value = "АВТ"
$ws1.Cells.Item(3,2).Value = "Автомобилист"
$ws1.Cells.Item(3,3).Value = "avtomobilist"
$ws1.Cells.Item(3,4).Value = "Денежкин Максим"
$ws1.Cells.Item(3,5).Value = "'24"
$ws1.Cells.Item(3,6).Value = "нападающий"
$ws1.Cells.Item(3,7).Value = "'25396"
$ws1.Cells.Item(3,8).Value = "1369_АВТ_денежкинмаксим"
$ws1.Cells.Item(3,9).Value = "injured_active"
$ws1.Cells.Item(3,10).Value = "https://www.khl.ru/clubs/avtomobilist/team/"
$ws1.Cells.Item(3,11).Value = "2025-12-13T03:01:21.677074+00:00"

$ws1.Cells.Item(4,1).Value = "АВТ"
$ws1.Cells.Item(4,2).Value = "Автомобилист"
$ws1.Cells.Item(4,3).Value = "avtomobilist"
$ws1.Cells.Item(4,4).Value = "Зборовский Сергей"
$ws1.Cells.Item(4,5).Value = "'2"
$ws1.Cells.Item(4,6).Value = "защитник"
$ws1.Cells.Item(4,7).Value = "'20989"
$ws1.Cells.Item(4,8).Value = "1369_АВТ_зборовскийсергей"
$ws1.Cells.Item(4,9).Value = "injured_active"
$ws1.Cells.Item(4,10).Value = "https://www.khl.ru/clubs/avtomobilist/team/"
$ws1.Cells.Item(4,11).Value = "2025-12-13T03:01:21.677094+00:00"

$ws1.Cells.Item(5,1).Value = "АДМ"
$ws1.Cells.Item(5,2).Value = "Адмирал"
$ws1.Cells.Item(5,3).Value = "admiral"
$ws1.Cells.Item(5,4).Value = "Дарьин Александр"
$ws1.Cells.Item(5,5).Value = "'41"
$ws1.Cells.Item(5,6).Value = "нападающий"
$ws1.Cells.Item(5,7).Value = "'26282"
$ws1.Cells.Item(5,8).Value = "1369_АДМ_дарьиналександр"
$ws1.Cells.Item(5,9).Value = "injured_active"
$ws1.Cells.Item(5,10).Value = "https://www.khl.ru/clubs/admiral/team/"
$ws1.Cells.Item(5,11).Value = "2025-12-13T03:01:24.583770+00:00"

$ws1.Cells.Item(6,1).Value = "АДМ"
$ws1.Cells.Item(6,2).Value = "Адмирал"
$ws1.Cells.Item(6,3).Value = "admiral"
$ws1.Cells.Item(6,4).Value = "Педан Руслан"
$ws1.Cells.Item(6,5).Value = "'53"
$ws1.Cells.Item(6,6).Value = "защитник"
$ws1.Cells.Item(6,7).Value = "'25353"
$ws1.Cells.Item(6,8).Value = "1369_АДМ_педанруслан"
$ws1.Cells.Item(6,9).Value = "injured_active"
$ws1.Cells.Item(6,10).Value = "https://www.khl.ru/clubs/admiral/team/"
$ws1.Cells.Item(6,11).Value = "2025-12-13T03:01:24.583797+00:00"

$ws1.Cells.Item(7,1).Value = "АДМ"
$ws1.Cells.Item(7,2).Value = "Адмирал"
$ws1.Cells.Item(7,3).Value = "admiral"
$ws1.Cells.Item(7,4).Value = "Цыба Арсений"
$ws1.Cells.Item(7,5).Value = "'1"
$ws1.Cells.Item(7,6).Value = "вратарь"
$ws1.Cells.Item(7,7).Value = "'33830"
$ws1.Cells.Item(7,8).Value = "1369_АДМ_цыбаарсений"
$ws1.Cells.Item(7,9).Value = "injured_active"
$ws1.Cells.Item(7,10).Value = "https://www.khl.ru/clubs/admiral/team/"
$ws1.Cells.Item(7,11).Value = "2025-12-13T03:01:24.583815+00:00"

$ws1.Cells.Item(8,1).Value = "АКБ"
$ws1.Cells.Item(8,2).Value = "Ак Барс"
$ws1.Cells.Item(8,3).Value = "ak_bars"
$ws1.Cells.Item(8,4).Value = "Яруллин Альберт"
$ws1.Cells.Item(8,5).Value = "'33"
$ws1.Cells.Item(8,6).Value = "защитник"
$ws1.Cells.Item(8,7).Value = "'16365"
$ws1.Cells.Item(8,8).Value = "1369_АКБ_яруллинальберт"
$ws1.Cells.Item(8,9).Value = "injured_active"
$ws1.Cells.Item(8,10).Value = "https://www.khl.ru/clubs/ak_bars/team/"
$ws1.Cells.Item(8,11).Value = "2025-12-13T03:01:26.935240+00:00"

$ws1.Cells.Item(9,1).Value = "АМР"
$ws1.Cells.Item(9,2).Value = "Амур"
$ws1.Cells.Item(9,3).Value = "amur"
$ws1.Cells.Item(9,4).Value = "Абросимов Роман"
$ws1.Cells.Item(9,5).Value = "'94"
$ws1.Cells.Item(9,6).Value = "защитник"
$ws1.Cells.Item(9,7).Value = "'17968"
$ws1.Cells.Item(9,8).Value = "1369_АМР_абросимовроман"
$ws1.Cells.Item(9,9).Value = "injured_active"
$ws1.Cells.Item(9,10).Value = "https://www.khl.ru/clubs/amur/team/"
$ws1.Cells.Item(9,11).Value = "2025-12-13T03:01:29.662491+00:00"

$ws1.Cells.Item(10,1).Value = "БАР"
$ws1.Cells.Item(10,2).Value = "Барыс"
$ws1.Cells.Item(10,3).Value = "barys"
$ws1.Cells.Item(10,4).Value = "Уотерспун Тайлер"
$ws1.Cells.Item(10,5).Value = "'26"
$ws1.Cells.Item(10,6).Value = "защитник"
$ws1.Cells.Item(10,7).Value = "'45769"
$ws1.Cells.Item(10,8).Value = "1369_БАР_уотерспунтайлер"
$ws1.Cells.Item(10,9).Value = "injured_active"
$ws1.Cells.Item(10,10).Value = "https://www.khl.ru/clubs/barys/team/"
$ws1.Cells.Item(10,11).Value = "2025-12-13T03:01:32.517913+00:00"

$ws1.Cells.Item(11,1).Value = "ДИН"
$ws1.Cells.Item(11,2).Value = "Динамо М"
$ws1.Cells.Item(11,3).Value = "dynamo_msk"
$ws1.Cells.Item(11,4).Value = "Готовец Кирилл"
$ws1.Cells.Item(11,5).Value = "'41"
$ws1.Cells.Item(11,6).Value = "защитник"
$ws1.Cells.Item(11,7).Value = "'16034"
$ws1.Cells.Item(11,8).Value = "1369_ДИН_готовецкирилл"
$ws1.Cells.Item(11,9).Value = "injured_active"
$ws1.Cells.Item(11,10).Value = "https://www.khl.ru/clubs/dynamo_msk/team/"
$ws1.Cells.Item(11,11).Value = "2025-12-13T03:01:35.373803+00:00"

$ws1.Cells.Item(12,1).Value = "ЛАД"
$ws1.Cells.Item(12,2).Value = "Лада"
$ws1.Cells.Item(12,3).Value = "lada"
$ws1.Cells.Item(12,4).Value = "Обидин Андрей"
$ws1.Cells.Item(12,5).Value = "'58"
$ws1.Cells.Item(12,6).Value = "нападающий"
$ws1.Cells.Item(12,7).Value = "'21310"
$ws1.Cells.Item(12,8).Value = "1369_ЛАД_обидинандрей"
$ws1.Cells.Item(12,9).Value = "injured_active"
$ws1.Cells.Item(12,10).Value = "https://www.khl.ru/clubs/lada/team/"
$ws1.Cells.Item(12,11).Value = "2025-12-13T03:01:40.416733+00:00"

$ws1.Cells.Item(13,1).Value = "ЛАД"
$ws1.Cells.Item(13,2).Value = "Лада"
$ws1.Cells.Item(13,3).Value = "lada"
$ws1.Cells.Item(13,4).Value = "Ожгихин Алексей"
$ws1.Cells.Item(13,5).Value = "'43"
$ws1.Cells.Item(13,6).Value = "нападающий"
$ws1.Cells.Item(13,7).Value = "'23021"
$ws1.Cells.Item(13,8).Value = "1369_ЛАД_ожгихиналексей"
$ws1.Cells.Item(13,9).Value = "injured_active"
$ws1.Cells.Item(13,10).Value = "https://www.khl.ru/clubs/lada/team/"
$ws1.Cells.Item(13,11).Value = "2025-12-13T03:01:40.416760+00:00"

$ws1.Cells.Item(14,1).Value = "ЛОК"
$ws1.Cells.Item(14,2).Value = "Локомотив"
$ws1.Cells.Item(14,3).Value = "lokomotiv"
$ws1.Cells.Item(14,4).Value = "Сергеев Андрей"
$ws1.Cells.Item(14,5).Value = "'99"
$ws1.Cells.Item(14,6).Value = "защитник"
$ws1.Cells.Item(14,7).Value = "'15416"
$ws1.Cells.Item(14,8).Value = "1369_ЛОК_сергеевандрей"
$ws1.Cells.Item(14,9).Value = "injured_active"
$ws1.Cells.Item(14,10).Value = "https://www.khl.ru/clubs/lokomotiv/team/"
$ws1.Cells.Item(14,11).Value = "2025-12-13T03:01:42.783776+00:00"

$ws1.Cells.Item(15,1).Value = "ММГ"
$ws1.Cells.Item(15,2).Value = "Металлург Мг"
$ws1.Cells.Item(15,3).Value = "metallurg_mg"
$ws1.Cells.Item(15,4).Value = "Сиряцкий Александр"
$ws1.Cells.Item(15,5).Value = "'74"
$ws1.Cells.Item(15,6).Value = "защитник"
$ws1.Cells.Item(15,7).Value = "'42458"
$ws1.Cells.Item(15,8).Value = "1369_ММГ_сиряцкийалександр"
$ws1.Cells.Item(15,9).Value = "injured_active"
$ws1.Cells.Item(15,10).Value = "https://www.khl.ru/clubs/metallurg_mg/team/"
$ws1.Cells.Item(15,11).Value = "2025-12-13T03:01:45.601631+00:00"

$ws1.Cells.Item(16,1).Value = "НХК"
$ws1.Cells.Item(16,2).Value = "Нефтехимик"
$ws1.Cells.Item(16,3).Value = "neftekhimik"
$ws1.Cells.Item(16,4).Value = "Хлыстов Никита"
$ws1.Cells.Item(16,5).Value = "'7"
$ws1.Cells.Item(16,6).Value = "защитник"
$ws1.Cells.Item(16,7).Value = "'17881"
$ws1.Cells.Item(16,8).Value = "1369_НХК_хлыстовникита"
$ws1.Cells.Item(16,9).Value = "injured_active"
$ws1.Cells.Item(16,10).Value = "https://www.khl.ru/clubs/neftekhimik/team/"
$ws1.Cells.Item(16,11).Value = "2025-12-13T03:01:48.281137+00:00"

$ws1.Cells.Item(17,1).Value = "СЕВ"
$ws1.Cells.Item(17,2).Value = "Северсталь"
$ws1.Cells.Item(17,3).Value = "severstal"
$ws1.Cells.Item(17,4).Value = "Ващенко Григорий"
$ws1.Cells.Item(17,5).Value = "'16"
$ws1.Cells.Item(17,6).Value = "защитник"
$ws1.Cells.Item(17,7).Value = "'14155"
$ws1.Cells.Item(17,8).Value = "1369_СЕВ_ващенкогригорий"
$ws1.Cells.Item(17,9).Value = "injured_active"
$ws1.Cells.Item(17,10).Value = "https://www.khl.ru/clubs/severstal/team/"
$ws1.Cells.Item(17,11).Value = "2025-12-13T03:01:50.583166+00:00"

$ws1.Cells.Item(18,1).Value = "СЕВ"
$ws1.Cells.Item(18,2).Value = "Северсталь"
$ws1.Cells.Item(18,3).Value = "severstal"
$ws1.Cells.Item(18,4).Value = "Смирнов Егор Д"
$ws1.Cells.Item(18,5).Value = "'47"
$ws1.Cells.Item(18,6).Value = "защитник"
$ws1.Cells.Item(18,7).Value = "'40906"
$ws1.Cells.Item(18,8).Value = "1369_СЕВ_смирновегорд"
$ws1.Cells.Item(18,9).Value = "injured_active"
$ws1.Cells.Item(18,10).Value = "https://www.khl.ru/clubs/severstal/team/"
$ws1.Cells.Item(18,11).Value = "2025-12-13T03:01:50.583196+00:00"

$ws1.Cells.Item(19,1).Value = "СЕВ"
$ws1.Cells.Item(19,2).Value = "Северсталь"
$ws1.Cells.Item(19,3).Value = "severstal"
$ws1.Cells.Item(19,4).Value = "Фомин Макар"
$ws1.Cells.Item(19,5).Value = "'77"
$ws1.Cells.Item(19,6).Value = "защитник"
$ws1.Cells.Item(19,7).Value = "'42087"
$ws1.Cells.Item(19,8).Value = "1369_СЕВ_фоминмакар"
$ws1.Cells.Item(19,9).Value = "injured_active"
$ws1.Cells.Item(19,10).Value = "https://www.khl.ru/clubs/severstal/team/"
$ws1.Cells.Item(19,11).Value = "2025-12-13T03:01:50.583212+00:00"

$ws1.Cells.Item(20,1).Value = "СЕВ"
$ws1.Cells.Item(20,2).Value = "Северсталь"
$ws1.Cells.Item(20,3).Value = "severstal"
$ws1.Cells.Item(20,4).Value = "Цицюра Владислав"
$ws1.Cells.Item(20,5).Value = "'10"
$ws1.Cells.Item(20,6).Value = "нападающий"
$ws1.Cells.Item(20,7).Value = "'23840"
$ws1.Cells.Item(20,8).Value = "1369_СЕВ_цицюравладислав"
$ws1.Cells.Item(20,9).Value = "injured_active"
$ws1.Cells.Item(20,10).Value = "https://www.khl.ru/clubs/severstal/team/"
$ws1.Cells.Item(20,11).Value = "2025-12-13T03:01:50.583228+00:00"

$ws1.Cells.Item(21,1).Value = "СЕВ"
$ws1.Cells.Item(21,2).Value = "Северсталь"
$ws1.Cells.Item(21,3).Value = "severstal"
$ws1.Cells.Item(21,4).Value = "Шостак Константин"
$ws1.Cells.Item(21,5).Value = "'1"
$ws1.Cells.Item(21,6).Value = "вратарь"
$ws1.Cells.Item(21,7).Value = "'27876"
$ws1.Cells.Item(21,8).Value = "1369_СЕВ_шостакконстантин"
$ws1.Cells.Item(21,9).Value = "injured_active"
$ws1.Cells.Item(21,10).Value = "https://www.khl.ru/clubs/severstal/team/"
$ws1.Cells.Item(21,11).Value = "2025-12-13T03:01:50.583243+00:00"

$ws1.Cells.Item(22,1).Value = "СИБ"
$ws1.Cells.Item(22,2).Value = "Сибирь"
$ws1.Cells.Item(22,3).Value = "sibir"
$ws1.Cells.Item(22,4).Value = "Аланов Егор"
$ws1.Cells.Item(22,5).Value = "'21"
$ws1.Cells.Item(22,6).Value = "защитник"
$ws1.Cells.Item(22,7).Value = "'26698"
$ws1.Cells.Item(22,8).Value = "1369_СИБ_алановегор"
$ws1.Cells.Item(22,9).Value = "injured_active"
$ws1.Cells.Item(22,10).Value = "https://www.khl.ru/clubs/sibir/team/"
$ws1.Cells.Item(22,11).Value = "2025-12-13T03:01:52.963855+00:00"

$ws1.Cells.Item(23,1).Value = "СИБ"
$ws1.Cells.Item(23,2).Value = "Сибирь"
$ws1.Cells.Item(23,3).Value = "sibir"
$ws1.Cells.Item(23,4).Value = "Першаков Александр"
$ws1.Cells.Item(23,5).Value = "'28"
$ws1.Cells.Item(23,6).Value = "нападающий"
$ws1.Cells.Item(23,7).Value = "'42479"
$ws1.Cells.Item(23,8).Value = "1369_СИБ_першаковалександр"
$ws1.Cells.Item(23,9).Value = "injured_active"
$ws1.Cells.Item(23,10).Value = "https://www.khl.ru/clubs/sibir/team/"
$ws1.Cells.Item(23,11).Value = "2025-12-13T03:01:52.963886+00:00"

$ws1.Cells.Item(24,1).Value = "СИБ"
$ws1.Cells.Item(24,2).Value = "Сибирь"
$ws1.Cells.Item(24,3).Value = "sibir"
$ws1.Cells.Item(24,4).Value = "Пьянов Валентин"
$ws1.Cells.Item(24,5).Value = "'45"
$ws1.Cells.Item(24,6).Value = "нападающий"
$ws1.Cells.Item(24,7).Value = "'16195"
$ws1.Cells.Item(24,8).Value = "1369_СИБ_пьяноввалентин"
$ws1.Cells.Item(24,9).Value = "injured_active"
$ws1.Cells.Item(24,10).Value = "https://www.khl.ru/clubs/sibir/team/"
$ws1.Cells.Item(24,11).Value = "2025-12-13T03:01:52.963903+00:00"

$ws1.Cells.Item(25,1).Value = "СКА"
$ws1.Cells.Item(25,2).Value = "СКА"
$ws1.Cells.Item(25,3).Value = "ska"
$ws1.Cells.Item(25,4).Value = "Зайцев Никита И"
$ws1.Cells.Item(25,5).Value = "'22"
$ws1.Cells.Item(25,6).Value = "защитник"
$ws1.Cells.Item(25,7).Value = "'16024"
$ws1.Cells.Item(25,8).Value = "1369_СКА_зайцевникитаи"
$ws1.Cells.Item(25,9).Value = "injured_active"
$ws1.Cells.Item(25,10).Value = "https://www.khl.ru/clubs/ska/team/"
$ws1.Cells.Item(25,11).Value = "2025-12-13T03:01:55.736817+00:00"

$ws1.Cells.Item(26,1).Value = "СКА"
$ws1.Cells.Item(26,2).Value = "СКА"
$ws1.Cells.Item(26,3).Value = "ska"
$ws1.Cells.Item(26,4).Value = "Зыков Валентин"
$ws1.Cells.Item(26,5).Value = "'90"
$ws1.Cells.Item(26,6).Value = "нападающий"
$ws1.Cells.Item(26,7).Value = "'17992"
$ws1.Cells.Item(26,8).Value = "1369_СКА_зыковвалентин"
$ws1.Cells.Item(26,9).Value = "injured_active"
$ws1.Cells.Item(26,10).Value = "https://www.khl.ru/clubs/ska/team/"
$ws1.Cells.Item(26,11).Value = "2025-12-13T03:01:55.736847+00:00"

$ws1.Cells.Item(27,1).Value = "СКА"
$ws1.Cells.Item(27,2).Value = "СКА"
$ws1.Cells.Item(27,3).Value = "ska"
$ws1.Cells.Item(27,4).Value = "Короткий Матвей"
$ws1.Cells.Item(27,5).Value = "'71"
$ws1.Cells.Item(27,6).Value = "нападающий"
$ws1.Cells.Item(27,7).Value = "'41566"
$ws1.Cells.Item(27,8).Value = "1369_СКА_короткийматвей"
$ws1.Cells.Item(27,9).Value = "injured_active"
$ws1.Cells.Item(27,10).Value = "https://www.khl.ru/clubs/ska/team/"
$ws1.Cells.Item(27,11).Value = "2025-12-13T03:01:55.736864+00:00"

$ws1.Cells.Item(28,1).Value = "СКА"
$ws1.Cells.Item(28,2).Value = "СКА"
$ws1.Cells.Item(28,3).Value = "ska"
$ws1.Cells.Item(28,4).Value = "Мёрфи Тревор"
$ws1.Cells.Item(28,5).Value = "'8"
$ws1.Cells.Item(28,6).Value = "защитник"
$ws1.Cells.Item(28,7).Value = "'34733"
$ws1.Cells.Item(28,8).Value = "1369_СКА_мерфитревор"
$ws1.Cells.Item(28,9).Value = "injured_active"
$ws1.Cells.Item(28,10).Value = "https://www.khl.ru/clubs/ska/team/"
$ws1.Cells.Item(28,11).Value = "2025-12-13T03:01:55.736881+00:00"

$ws1.Cells.Item(29,1).Value = "СКА"
$ws1.Cells.Item(29,2).Value = "СКА"
$ws1.Cells.Item(29,3).Value = "ska"
$ws1.Cells.Item(29,4).Value = "Сапего Сергей"
$ws1.Cells.Item(29,5).Value = "'9"
$ws1.Cells.Item(29,6).Value = "защитник"
$ws1.Cells.Item(29,7).Value = "'39875"
$ws1.Cells.Item(29,8).Value = "1369_СКА_сапегосергей"
$ws1.Cells.Item(29,9).Value = "injured_active"
$ws1.Cells.Item(29,10).Value = "https://www.khl.ru/clubs/ska/team/"
$ws1.Cells.Item(29,11).Value = "2025-12-13T03:01:55.736897+00:00"

$ws1.Cells.Item(30,1).Value = "СОЧ"
$ws1.Cells.Item(30,2).Value = "ХК Сочи"
$ws1.Cells.Item(30,3).Value = "hc_sochi"
$ws1.Cells.Item(30,4).Value = "Бикмуллин Рафаэль"
$ws1.Cells.Item(30,5).Value = "'24"
$ws1.Cells.Item(30,6).Value = "нападающий"
$ws1.Cells.Item(30,7).Value = "'22424"
$ws1.Cells.Item(30,8).Value = "1369_СОЧ_бикмуллинрафаэль"
$ws1.Cells.Item(30,9).Value = "injured_active"
$ws1.Cells.Item(30,10).Value = "https://www.khl.ru/clubs/hc_sochi/team/"
$ws1.Cells.Item(30,11).Value = "2025-12-13T03:01:58.999730+00:00"

$ws1.Cells.Item(31,1).Value = "СОЧ"
$ws1.Cells.Item(31,2).Value = "ХК Сочи"
$ws1.Cells.Item(31,3).Value = "hc_sochi"
$ws1.Cells.Item(31,4).Value = "Венгрыжановский Денис"
$ws1.Cells.Item(31,5).Value = "'9"
$ws1.Cells.Item(31,6).Value = "нападающий"
$ws1.Cells.Item(31,7).Value = "'31892"
$ws1.Cells.Item(31,8).Value = "1369_СОЧ_венгрыжановскийденис"
$ws1.Cells.Item(31,9).Value = "injured_active"
$ws1.Cells.Item(31,10).Value = "https://www.khl.ru/clubs/hc_sochi/team/"
$ws1.Cells.Item(31,11).Value = "2025-12-13T03:01:58.999763+00:00"

$ws1.Cells.Item(32,1).Value = "СОЧ"
$ws1.Cells.Item(32,2).Value = "ХК Сочи"
$ws1.Cells.Item(32,3).Value = "hc_sochi"
$ws1.Cells.Item(32,4).Value = "Хёфенмайер Ноэль"
$ws1.Cells.Item(32,5).Value = "'22"
$ws1.Cells.Item(32,6).Value = "защитник"
$ws1.Cells.Item(32,7).Value = "'44847"
$ws1.Cells.Item(32,8).Value = "1369_СОЧ_хефенмайерноэль"
$ws1.Cells.Item(32,9).Value = "injured_active"
$ws1.Cells.Item(32,10).Value = "https://www.khl.ru/clubs/hc_sochi/team/"
$ws1.Cells.Item(32,11).Value = "2025-12-13T03:01:58.999781+00:00"

$ws1.Cells.Item(33,1).Value = "СПР"
$ws1.Cells.Item(33,2).Value = "Спартак"
$ws1.Cells.Item(33,3).Value = "spartak"
$ws1.Cells.Item(33,4).Value = "Вишневский Дмитрий"
$ws1.Cells.Item(33,5).Value = "'55"
$ws1.Cells.Item(33,6).Value = "защитник"
$ws1.Cells.Item(33,7).Value = "'15299"
$ws1.Cells.Item(33,8).Value = "1369_СПР_вишневскийдмитрий"
$ws1.Cells.Item(33,9).Value = "injured_active"
$ws1.Cells.Item(33,10).Value = "https://www.khl.ru/clubs/spartak/team/"
$ws1.Cells.Item(33,11).Value = "2025-12-13T03:02:01.834154+00:00"

$ws1.Cells.Item(34,1).Value = "СПР"
$ws1.Cells.Item(34,2).Value = "Спартак"
$ws1.Cells.Item(34,3).Value = "spartak"
$ws1.Cells.Item(34,4).Value = "Коростелёв Никита"
$ws1.Cells.Item(34,5).Value = "'35"
$ws1.Cells.Item(34,6).Value = "нападающий"
$ws1.Cells.Item(34,7).Value = "'22149"
$ws1.Cells.Item(34,8).Value = "1369_СПР_коростелевникита"
$ws1.Cells.Item(34,9).Value = "injured_active"
$ws1.Cells.Item(34,10).Value = "https://www.khl.ru/clubs/spartak/team/"
$ws1.Cells.Item(34,11).Value = "2025-12-13T03:02:01.834182+00:00"

$ws1.Cells.Item(35,1).Value = "СПР"
$ws1.Cells.Item(35,2).Value = "Спартак"
$ws1.Cells.Item(35,3).Value = "spartak"
$ws1.Cells.Item(35,4).Value = "Порядин Павел"
$ws1.Cells.Item(35,5).Value = "'24"
$ws1.Cells.Item(35,6).Value = "нападающий"
$ws1.Cells.Item(35,7).Value = "'19258"
$ws1.Cells.Item(35,8).Value = "1369_СПР_порядинпавел"
$ws1.Cells.Item(35,9).Value = "injured_active"
$ws1.Cells.Item(35,10).Value = "https://www.khl.ru/clubs/spartak/team/"
$ws1.Cells.Item(35,11).Value = "2025-12-13T03:02:01.834199+00:00"

$ws1.Cells.Item(36,1).Value = "СЮЛ"
$ws1.Cells.Item(36,2).Value = "Салават Юлаев"
$ws1.Cells.Item(36,3).Value = "salavat_yulaev"
$ws1.Cells.Item(36,4).Value = "Алалыкин Данил"
$ws1.Cells.Item(36,5).Value = "'61"
$ws1.Cells.Item(36,6).Value = "нападающий"
$ws1.Cells.Item(36,7).Value = "'34493"
$ws1.Cells.Item(36,8).Value = "1369_СЮЛ_алалыкинданил"
$ws1.Cells.Item(36,9).Value = "injured_active"
$ws1.Cells.Item(36,10).Value = "https://www.khl.ru/clubs/salavat_yulaev/team/"
$ws1.Cells.Item(36,11).Value = "2025-12-13T03:02:04.564870+00:00"

$ws1.Cells.Item(37,1).Value = "ТОР"
$ws1.Cells.Item(37,2).Value = "Торпедо"
$ws1.Cells.Item(37,3).Value = "torpedo"
$ws1.Cells.Item(37,4).Value = "Воронин Кирилл"
$ws1.Cells.Item(37,5).Value = "'41"
$ws1.Cells.Item(37,6).Value = "нападающий"
$ws1.Cells.Item(37,7).Value = "'17354"
$ws1.Cells.Item(37,8).Value = "1369_ТОР_воронинкирилл"
$ws1.Cells.Item(37,9).Value = "injured_active"
$ws1.Cells.Item(37,10).Value = "https://www.khl.ru/clubs/torpedo/team/"
$ws1.Cells.Item(37,11).Value = "2025-12-13T03:02:07.382345+00:00"

$ws1.Cells.Item(38,1).Value = "ТОР"
$ws1.Cells.Item(38,2).Value = "Торпедо"
$ws1.Cells.Item(38,3).Value = "torpedo"
$ws1.Cells.Item(38,4).Value = "Кручинин Алексей"
$ws1.Cells.Item(38,5).Value = "'78"
$ws1.Cells.Item(38,6).Value = "нападающий"
$ws1.Cells.Item(38,7).Value = "'16355"
$ws1.Cells.Item(38,8).Value = "1369_ТОР_кручининалексей"
$ws1.Cells.Item(38,9).Value = "injured_active"
$ws1.Cells.Item(38,10).Value = "https://www.khl.ru/clubs/torpedo/team/"
$ws1.Cells.Item(38,11).Value = "2025-12-13T03:02:07.382375+00:00"

$ws1.Cells.Item(39,1).Value = "ТОР"
$ws1.Cells.Item(39,2).Value = "Торпедо"
$ws1.Cells.Item(39,3).Value = "torpedo"
$ws1.Cells.Item(39,4).Value = "Принс Шэйн"
$ws1.Cells.Item(39,5).Value = "'18"
$ws1.Cells.Item(39,6).Value = "нападающий"
$ws1.Cells.Item(39,7).Value = "'19060"
$ws1.Cells.Item(39,8).Value = "1369_ТОР_принсшэйн"
$ws1.Cells.Item(39,9).Value = "injured_active"
$ws1.Cells.Item(39,10).Value = "https://www.khl.ru/clubs/torpedo/team/"
$ws1.Cells.Item(39,11).Value = "2025-12-13T03:02:07.382392+00:00"

$ws1.Cells.Item(40,1).Value = "ТРК"
$ws1.Cells.Item(40,2).Value = "Трактор"
$ws1.Cells.Item(40,3).Value = "traktor"
$ws1.Cells.Item(40,4).Value = "Светлаков Андрей"
$ws1.Cells.Item(40,5).Value = "'87"
$ws1.Cells.Item(40,6).Value = "нападающий"
$ws1.Cells.Item(40,7).Value = "'19218"
$ws1.Cells.Item(40,8).Value = "1369_ТРК_светлаковандрей"
$ws1.Cells.Item(40,9).Value = "injured_active"
$ws1.Cells.Item(40,10).Value = "https://www.khl.ru/clubs/traktor/team/"
$ws1.Cells.Item(40,11).Value = "2025-12-13T03:02:09.800178+00:00"

$ws1.Cells.Item(41,1).Value = "ЦСК"
$ws1.Cells.Item(41,2).Value = "ЦСКА"
$ws1.Cells.Item(41,3).Value = "cska"
$ws1.Cells.Item(41,4).Value = "Бучельников Дмитрий"
$ws1.Cells.Item(41,5).Value = "'72"
$ws1.Cells.Item(41,6).Value = "нападающий"
$ws1.Cells.Item(41,7).Value = "'39102"
$ws1.Cells.Item(41,8).Value = "1369_ЦСК_бучельниковдмитрий"
$ws1.Cells.Item(41,9).Value = "injured_active"
$ws1.Cells.Item(41,10).Value = "https://www.khl.ru/clubs/cska/team/"
$ws1.Cells.Item(41,11).Value = "2025-12-13T03:02:12.690670+00:00"

$ws1.Cells.Item(42,1).Value = "ЦСК"
$ws1.Cells.Item(42,2).Value = "ЦСКА"
$ws1.Cells.Item(42,3).Value = "cska"
$ws1.Cells.Item(42,4).Value = "Моисеев Данила"
$ws1.Cells.Item(42,5).Value = "'93"
$ws1.Cells.Item(42,6).Value = "нападающий"
$ws1.Cells.Item(42,7).Value = "'23931"
$ws1.Cells.Item(42,8).Value = "1369_ЦСК_моисеевданила"
$ws1.Cells.Item(42,9).Value = "injured_active"
$ws1.Cells.Item(42,10).Value = "https://www.khl.ru/clubs/cska/team/"
$ws1.Cells.Item(42,11).Value = "2025-12-13T03:02:12.690699+00:00"

$ws1.Cells.Item(43,1).Value = "ШДР"
$ws1.Cells.Item(43,2).Value = "Драконы"
$ws1.Cells.Item(43,3).Value = "kunlun"
$ws1.Cells.Item(43,4).Value = "Бишофф Джейк"
$ws1.Cells.Item(43,5).Value = "'28"
$ws1.Cells.Item(43,6).Value = "защитник"
$ws1.Cells.Item(43,7).Value = "'45490"
$ws1.Cells.Item(43,8).Value = "1369_ШДР_бишоффджейк"
$ws1.Cells.Item(43,9).Value = "injured_active"
$ws1.Cells.Item(43,10).Value = "https://www.khl.ru/clubs/kunlun/team/"
$ws1.Cells.Item(43,11).Value = "2025-12-13T03:02:15.420781+00:00"

$ws1.Cells.Item(44,1).Value = "ШДР"
$ws1.Cells.Item(44,2).Value = "Драконы"
$ws1.Cells.Item(44,3).Value = "kunlun"
$ws1.Cells.Item(44,4).Value = "Гроло Жереми"
$ws1.Cells.Item(44,5).Value = "'75"
$ws1.Cells.Item(44,6).Value = "защитник"
$ws1.Cells.Item(44,7).Value = "'45343"
$ws1.Cells.Item(44,8).Value = "1369_ШДР_гроложереми"
$ws1.Cells.Item(44,9).Value = "injured_active"
$ws1.Cells.Item(44,10).Value = "https://www.khl.ru/clubs/kunlun/team/"
$ws1.Cells.Item(44,11).Value = "2025-12-13T03:02:15.420808+00:00"

# ---- Sheet: returned ----
$ws2 = $wb.Worksheets.Item("returned")

$ws2.Cells.Item(2,1).Value = "ТРК"
$ws2.Cells.Item(2,2).Value = "Трактор"
$ws2.Cells.Item(2,3).Value = "Мыльников Сергей И"
$ws2.Cells.Item(2,4).Value = "1369_ТРК_мыльниковсергейи"
$ws2.Cells.Item(2,5).Value = "RETURN"
$ws2.Cells.Item(2,6).Value = "2025-12-13T11:02:15.927763+08:00"
$ws2.Cells.Item(2,7).Value = "'2025-12-13"

# ---- Sheet: new_injured ----
$ws3 = $wb.Worksheets.Item("new_injured")

$ws3.Cells.Item(2,1).Value = "СЕВ"
$ws3.Cells.Item(2,2).Value = "Северсталь"
$ws3.Cells.Item(2,3).Value = "Ващенко Григорий"
$ws3.Cells.Item(2,4).Value = "1369_СЕВ_ващенкогригорий"
$ws3.Cells.Item(2,5).Value = "INJURED_NEW"
$ws3.Cells.Item(2,6).Value = "2025-12-13T11:02:15.927763+08:00"
$ws3.Cells.Item(2,7).Value = "'2025-12-13"

$ws3.Cells.Item(3,1).Value = "СЕВ"
$ws3.Cells.Item(3,2).Value = "Северсталь"
$ws3.Cells.Item(3,3).Value = "Смирнов Егор Д"
$ws3.Cells.Item(3,4).Value = "1369_СЕВ_смирновегорд"
$ws3.Cells.Item(3,5).Value = "INJURED_NEW"
$ws3.Cells.Item(3,6).Value = "2025-12-13T11:02:15.927763+08:00"
$ws3.Cells.Item(3,7).Value = "'2025-12-13"

$ws3.Cells.Item(4,1).Value = "СЕВ"
$ws3.Cells.Item(4,2).Value = "Северсталь"
$ws3.Cells.Item(4,3).Value = "Шостак Константин"
$ws3.Cells.Item(4,4).Value = "1369_СЕВ_шостакконстантин"
$ws3.Cells.Item(4,5).Value = "INJURED_NEW"
$ws3.Cells.Item(4,6).Value = "2025-12-13T11:02:15.927763+08:00"
$ws3.Cells.Item(4,7).Value = "'2025-12-13"

$ws3.Cells.Item(5,1).Value = "СПР"
$ws3.Cells.Item(5,2).Value = "Спартак"
$ws3.Cells.Item(5,3).Value = "Порядин Павел"
$ws3.Cells.Item(5,4).Value = "1369_СПР_порядинпавел"
$ws3.Cells.Item(5,5).Value = "INJURED_NEW"
$ws3.Cells.Item(5,6).Value = "2025-12-13T11:02:15.927763+08:00"
$ws3.Cells.Item(5,7).Value = "'2025-12-13"

